$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.23"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "24.03"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.237"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05787"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.501"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8612"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1368"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06958"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03194"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02879"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09379"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.755"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001515"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04717"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006010"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006272"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001236"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.004546"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00006900"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.510"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.122"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3190"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1357"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03652"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006280"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1053"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003400"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008032"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005275"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.3320"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002350"
